$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 8-12 (the underlying observation records got
# reshuffled between rows, per the source diff). "I"/"J" are only listed
# when the diff actually shows a change for that row; the other rows'
# (already blank) Antal/Enhet cells are left untouched.
$rowData = @{
    8  = @{ A = 111702420; B = 90709; D = "NT"; E = 5448; F = "Svartvit taggsvamp"; G = "Phellodon connatus"; H = "(Schultz) nom.prov"; I = "'1"; J = "fruktkroppar"; P = "Kyrkogården (Kyrkogården), Nrk"; Q = 517086.1792710476; R = 6574909.900584662 }
    9  = @{ A = 111702393; B = 89183;  D = "LC"; E = 3215; F = "Rödgul trumpetsvamp"; G = "Craterellus lutescens"; H = "(Fr.) Fr."; P = "Kyrkogården (Kyrkogården), Nrk"; Q = 517070.2129045375; R = 6574934.844418272 }
    10 = @{ A = 111702486; B = 90678;  D = "LC"; E = 4366; F = "Skarp dropptaggsvamp"; G = "Hydnellum peckii"; H = "Banker"; P = "Kyrkogården (Kyrkogården), Nrk"; Q = 517080.8398438052; R = 6574959.907818918 }
    11 = @{ A = 111702400; B = 90687;  D = "LC"; E = 5964; F = "Fjällig taggsvamp s.str."; G = "Sarcodon imbricatus s.str."; H = "(L.:Fr.) P.Karst."; I = "CLEAR"; J = "CLEAR"; P = "Kyrkogården (Kyrkogården), Nrk"; Q = 517073.2951468225; R = 6574931.795150192 }
    12 = @{ A = 111702506; B = 90687;  D = "LC"; E = 5964; F = "Fjällig taggsvamp s.str."; G = "Sarcodon imbricatus s.str."; H = "(L.:Fr.) P.Karst."; P = "Kyrkogården, Nrk"; Q = 517093.6249861007; R = 6574959.965416327 }
}

foreach ($r in $rowData.Keys) {
    $d = $rowData[$r]
    $ws.Cells.Item($r, 1).Value = $d.A
    $ws.Cells.Item($r, 2).Value = $d.B
    $ws.Cells.Item($r, 4).Value = $d.D
    $ws.Cells.Item($r, 5).Value = $d.E
    $ws.Cells.Item($r, 6).Value = $d.F
    $ws.Cells.Item($r, 7).Value = $d.G
    $ws.Cells.Item($r, 8).Value = $d.H

    if ($d.ContainsKey("I")) {
        if ($d.I -eq "CLEAR") {
            $ws.Cells.Item($r, 9).ClearContents()
        } else {
            $ws.Cells.Item($r, 9).Value = $d.I
        }
    }

    if ($d.ContainsKey("J")) {
        if ($d.J -eq "CLEAR") {
            $ws.Cells.Item($r, 10).ClearContents()
        } else {
            $ws.Cells.Item($r, 10).Value = $d.J
        }
    }

    $ws.Cells.Item($r, 16).Value = $d.P
    $ws.Cells.Item($r, 17).Value = $d.Q
    $ws.Cells.Item($r, 18).Value = $d.R
}
